$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in column H, matching the style used by the
# other header cells (bold font, thin border, centered/top aligned).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the new Save values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
